$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the rows that are no longer part of the table ---
# "More specific truck entries (do not implement these for now)" block (old rows 7-14)
$ws.Rows("7:14").Delete()
# "Transport, aircraft, freight" (old row 2)
$ws.Rows("2").Delete()

# --- Re-label the impact-category headers (row 1) ---
$ws.Range("B1").Value = "acp"
$ws.Range("C1").Value = "eup"
$ws.Range("D1").Value = "GWPf"
$ws.Range("E1").Value = "GWPb"
$ws.Range("F1").Value = "GWP-LULUC"
$ws.Range("G1").Value = "odp"
$ws.Range("H1").Value = "smg"

# --- Add the new formatted (empty) rows 15-21 ---
$first = $ws.Range("C15")
$first.Font.Name = "Consolas"
$first.Font.Color = 7901646
$first.VerticalAlignment = -4108

$rest = $ws.Range("C16:C21")
$rest.Font.Name = "Consolas"
$rest.Font.Color = 13421772
$rest.VerticalAlignment = -4108

# --- Page setup / view tweaks ---
$ws.PageSetup.Orientation = 1
[void]$ws.Range("F2").Select()
